$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the 5 "good" urls that were sitting in column B (B2:B6) down into
# column A, appended after the existing list (A11:A15) - this carries over
# both their values and their formatting, replacing the old entries that
# used to occupy A11:A13.
$ws.Range("B2:B6").Copy($ws.Range("A11:A15"))

# Column B is no longer used.
# B1 (the "bad_url" header) only loses its text, formatting is kept.
$ws.Range("B1").ClearContents()
# B2:B7 (the url values) are removed completely, formatting and all.
$ws.Range("B2:B7").Clear()

# Update the column A header text to reflect the new single-column layout.
$ws.Range("A1").Value = "recording_url"

# Move the active selection as recorded by the author.
$null = $ws.Range("C4").Select()
